# Word COM-interop script implementing the commit:
# "Three sets of 1000 word documents with three font types"
#
# Changes applied to this document (0627.docx):
#  1. Rename font "TimesNewToman" -> "Times New Roman" throughout.
#  2. Replace title / author / email placeholder text.
#  3. Replace the body paragraph and summary paragraph sentences with
#     new "Biology" themed content.
#  4. Remove the final sentence of the summary paragraph (it is dropped).
#  5. Append one extra empty paragraph at the end of the document.

$d = $word.ActiveDocument
$find = $d.Content.Find

# ---------------------------------------------------------------------
# 1. Fix the misspelled font name across the whole document.
#    (Font names live in run-properties, not in visible text, so this
#    must be done through the Font object rather than Find/Replace.)
# ---------------------------------------------------------------------
$fullRange = $d.Range(0, $d.Content.End)
$fullRange.Font.Name = "Times New Roman"

# ---------------------------------------------------------------------
# 2. Title
# ---------------------------------------------------------------------
$find.Execute("The Heart's Symphony: Rhythm of Life", $true, $false, $false, $false, $false, $true, 1, $false, "The Enchanting Realm of Biology: Unveiling the Secrets of Life", 2) | Out-Null

# ---------------------------------------------------------------------
# 3. Author name
# ---------------------------------------------------------------------
$find.Execute("Juliet Elizabeth", $true, $false, $false, $false, $false, $true, 1, $false, "Alex Kingston", 2) | Out-Null

# ---------------------------------------------------------------------
# 4. Email line: collapse the five runs (juliet / . / elizabeth@remedyhealth
#    / . / net) into a single run reading "valid_email_address". We find
#    the paragraph that currently holds the e-mail address and replace
#    its whole range (minus the paragraph mark) directly, which merges
#    the runs the same way Word does when you retype a selection.
# ---------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -eq "juliet.elizabeth@remedyhealth.net") {
        $emailRange = $d.Range($para.Range.Start, $para.Range.Start + 34)
        $emailRange.Text = "valid_email_address"
        break
    }
}

# ---------------------------------------------------------------------
# 5. Body paragraph sentences (first big paragraph).
# ---------------------------------------------------------------------
$find.Execute("Within the intricate chambers of the human body, a remarkable orchestra conducts the symphony of life", $true, $false, $false, $false, $false, $true, 1, $false, "Biology, the profound science that unravels the enigmatic tapestry of life, holds an intrinsic fascination for inquisitive minds", 2) | Out-Null
$find.Execute(" The heart, a tireless maestro, orchestrates a harmonious rhythm, pumping blood through an elaborate network of vessels, nourishing cells, and sustaining existence", $true, $false, $false, $false, $false, $true, 1, $false, " Embarking on this extraordinary journey, we delve into the intricate mechanisms that govern the diversity and complexity of living organisms", 2) | Out-Null
$find.Execute(" This vital organ, a symbol of resilience and unwavering dedication, plays a pivotal role in maintaining our physical and emotional well-being", $true, $false, $false, $false, $false, $true, 1, $false, " From the tiniest cellular components to the grand symphony of ecosystems, biology unveils the awe-inspiring beauty and interconnectedness of all life", 2) | Out-Null
$find.Execute("The heart, with its intricate structure and unwavering rhythm, has captivated the imagination of poets, philosophers, and scientists throughout history", $true, $false, $false, $false, $false, $true, 1, $false, "In the microscopic realm, biology unravels the symphony of life's fundamental unit, the cell", 2) | Out-Null
$find.Execute(" Celebrated in art, literature, and music, the heart's symbolic significance transcends its physiological function", $true, $false, $false, $false, $false, $true, 1, $false, " Delving into the depths of cellular structure, we discover the intricate interplay of organelles, each with its unique function, orchestrating the harmonious functioning of the whole", 2) | Out-Null
$find.Execute(" It represents love, passion, courage, and the very essence of human experience", $true, $false, $false, $false, $false, $true, 1, $false, " These microscopic marvels, the building blocks of all living things, reveal the enigmatic dance of molecules, the ceaseless flow of energy, and the intricate genetic code that holds the blueprint for life", 2) | Out-Null
$find.Execute("The medical and scientific exploration of the heart has yielded profound insights into its intricate workings", $true, $false, $false, $false, $false, $true, 1, $false, "Venturing beyond the cellular realm, biology unveils the staggering diversity of life forms that inhabit our planet", 2) | Out-Null
$find.Execute(" Advanced imaging techniques, such as echocardiography and cardiac MRI, allow physicians to visualize the heart in motion, revealing its chambers, valves, and blood flow patterns", $true, $false, $false, $false, $false, $true, 1, $false, " From the towering sequoia trees to the microscopic plankton, the Earth brims with a mesmerizing array of organisms, each exquisitely adapted to its specific niche", 2) | Out-Null
$find.Execute(" Electrocardiography records the heart's electrical activity, enabling the detection of abnormalities and guiding diagnosis and treatment", $true, $false, $false, $false, $false, $true, 1, $false, " This symphony of life, a testament to evolution's creative genius, highlights the remarkable resilience and adaptability of living organisms in the face of a dynamic and ever-changing environment", 2) | Out-Null

# ---------------------------------------------------------------------
# 6. Summary paragraph sentences.
# ---------------------------------------------------------------------
$find.Execute("The heart, a vital organ of the human body, plays a crucial role in maintaining physical and emotional well-being", $true, $false, $false, $false, $false, $true, 1, $false, "Biology, a science of awe and wonder, delves into the intricacies of life's fundamental unit, the cell, exploring the harmonious interplay of its organelles", 2) | Out-Null
$find.Execute(" Symbolized by love, passion, and courage, the heart has captured the imagination of poets, philosophers, and scientists throughout history", $true, $false, $false, $false, $false, $true, 1, $false, " It unravels the captivating diversity of life forms that inhabit our planet, revealing the remarkable resilience and adaptability of organisms in a dynamic environment", 2) | Out-Null
$find.Execute(" Medical advancements have enabled profound insights into the heart's intricate workings, aiding diagnosis and treatment of cardiac conditions", $true, $false, $false, $false, $false, $true, 1, $false, " Through the lens of biology, we gain a profound appreciation for the intricate beauty and interconnectedness of all life, fostering a sense of stewardship and responsibility toward the natural world", 2) | Out-Null

# ---------------------------------------------------------------------
# 7. Drop the final sentence of the summary paragraph entirely
#    (" As research continues ... existence." plus its trailing period
#    run) -- simply replace it with nothing.
# ---------------------------------------------------------------------
$find.Execute(" As research continues to unveil the heart's secrets, we marvel at its tireless rhythm, a testament to the wonders of human existence.", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# ---------------------------------------------------------------------
# 8. Append a new, empty paragraph at the very end of the document.
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$endRange = $d.Range($lastPara.Range.End, $lastPara.Range.End)
$endRange.InsertParagraphAfter() | Out-Null
